# Auto-generated edit script applying cell value changes to H:N columns
# across multiple worksheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 78.75
$ws.Range("I2").Value = 78.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 78.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 34.25
$ws.Range("N2").ClearContents()
$ws.Range("H8").Value = 156.85715
$ws.Range("I8").Value = 179.6
$ws.Range("K8").Value = 538.8
$ws.Range("M8").Value = -399.8
$ws.Range("H11").Value = 974.2857
$ws.Range("I11").Value = 974.2857
$ws.Range("K11").Value = 974.2857
$ws.Range("M11").Value = -834.2857
$ws.Range("H21").Value = 21999.5
$ws.Range("I21").Value = 20000
$ws.Range("J21").Value = 23999
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 23999
$ws.Range("M21").Value = -19532
$ws.Range("N21").Value = -24935
$ws.Range("H23").Value = 21999.5
$ws.Range("I23").Value = 20000
$ws.Range("J23").Value = 23999
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 23999
$ws.Range("M23").Value = -19766
$ws.Range("N23").Value = -24467
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 100
$ws.Range("K29").Value = 300
$ws.Range("M29").Value = -19
$ws.Range("H38").Value = 2240.2222
$ws.Range("J38").Value = 6333.3335
$ws.Range("L38").Value = 19000.0005
$ws.Range("N38").Value = -19744.0005
$ws.Range("H40").Value = 61495.133
$ws.Range("J40").Value = 3930.8086
$ws.Range("L40").Value = 3930.8086
$ws.Range("N40").Value = -4280.8086
$ws.Range("H42").Value = 123.916664
$ws.Range("I42").Value = 137.125
$ws.Range("J42").Value = 97.5
$ws.Range("K42").Value = 411.375
$ws.Range("L42").Value = 292.5
$ws.Range("M42").Value = -181.375
$ws.Range("N42").Value = -752.5
$ws.Range("H43").Value = 2194538
$ws.Range("I43").Value = 4924961
$ws.Range("K43").Value = 4924961
$ws.Range("M43").Value = -4924892
$ws.Range("H58").Value = 2655.8
$ws.Range("I58").Value = 3526.3333
$ws.Range("J58").Value = 1350
$ws.Range("K58").Value = 10578.9999
$ws.Range("L58").Value = 4050
$ws.Range("M58").Value = -10428.9999
$ws.Range("N58").Value = -4350
$ws.Range("H62").Value = 8784.532999999999
$ws.Range("I62").Value = 7898.1665
$ws.Range("K62").Value = 7898.1665
$ws.Range("M62").Value = -7274.1665
$ws.Range("H65").Value = 8784.532999999999
$ws.Range("I65").Value = 7898.1665
$ws.Range("K65").Value = 39490.8325
$ws.Range("M65").Value = -36370.8325
$ws.Range("H74").Value = 5997
$ws.Range("I74").Value = 5209.7144
$ws.Range("K74").Value = 5209.7144
$ws.Range("M74").Value = -4273.7144
$ws.Range("H77").Value = 5997
$ws.Range("I77").Value = 5209.7144
$ws.Range("K77").Value = 26048.572
$ws.Range("M77").Value = -21368.572
$ws.Range("H116").Value = 76500
$ws.Range("I116").Value = 76500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 76500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -73058
$ws.Range("N116").ClearContents()
$ws.Range("H135").Value = 1318.6
$ws.Range("I135").Value = 1318.6
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 11867.4
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9332.4
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 41039.57
$ws.Range("I137").Value = 165959.8
$ws.Range("J137").Value = 2002
$ws.Range("K137").Value = 497879.4
$ws.Range("L137").Value = 6006
$ws.Range("M137").Value = -495329.4
$ws.Range("N137").Value = -11106
$ws.Range("H138").Value = 3817.975
$ws.Range("I138").Value = 1479.2632
$ws.Range("J138").Value = 5933.952
$ws.Range("K138").Value = 4437.7896
$ws.Range("L138").Value = 17801.856
$ws.Range("M138").Value = 702.2103999999999
$ws.Range("N138").Value = -28081.856
$ws.Range("H141").Value = 4404.75
$ws.Range("I141").Value = 4404.75
$ws.Range("K141").Value = 13214.25
$ws.Range("M141").Value = -8034.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2078.889
$ws.Range("I2").Value = 1393.3334
$ws.Range("J2").Value = 3450
$ws.Range("K2").Value = 1393.3334
$ws.Range("L2").Value = 3450
$ws.Range("M2").Value = -1280.3334
$ws.Range("N2").Value = -3676
$ws.Range("H14").Value = 16637167
$ws.Range("I14").Value = 16637167
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 16637167
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -16636992
$ws.Range("N14").ClearContents()
$ws.Range("H34").Value = 136605
$ws.Range("I34").Value = 51008.332
$ws.Range("J34").Value = 265000
$ws.Range("K34").Value = 51008.332
$ws.Range("L34").Value = 265000
$ws.Range("M34").Value = -50737.332
$ws.Range("N34").Value = -265542
$ws.Range("H40").Value = 14000
$ws.Range("J40").Value = 14000
$ws.Range("L40").Value = 14000
$ws.Range("N40").Value = -14352
$ws.Range("H44").Value = 69995
$ws.Range("J44").Value = 69995
$ws.Range("L44").Value = 69995
$ws.Range("N44").Value = -70971
$ws.Range("H61").Value = 3008.611
$ws.Range("I61").Value = 2697.25
$ws.Range("K61").Value = 2697.25
$ws.Range("M61").Value = -2485.25
$ws.Range("H63").Value = 159088.89
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 284360
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 284360
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -285732
$ws.Range("H66").Value = 159088.89
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 284360
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 1421800
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -1428664
$ws.Range("H74").Value = 2237.4707
$ws.Range("I74").Value = 2002.5333
$ws.Range("J74").Value = 3999.5
$ws.Range("K74").Value = 2002.5333
$ws.Range("L74").Value = 3999.5
$ws.Range("M74").Value = -1128.5333
$ws.Range("N74").Value = -5747.5
$ws.Range("H77").Value = 2237.4707
$ws.Range("I77").Value = 2002.5333
$ws.Range("J77").Value = 3999.5
$ws.Range("K77").Value = 10012.6665
$ws.Range("L77").Value = 19997.5
$ws.Range("M77").Value = -5644.666499999999
$ws.Range("N77").Value = -28733.5
$ws.Range("H97").Value = 1191.7188
$ws.Range("I97").Value = 1038.88
$ws.Range("K97").Value = 1038.88
$ws.Range("M97").Value = -542.8800000000001
$ws.Range("H102").Value = 1793.9
$ws.Range("I102").Value = 1304.875
$ws.Range("K102").Value = 1304.875
$ws.Range("M102").Value = 317.125
$ws.Range("H114").Value = 47000
$ws.Range("J114").Value = 47000
$ws.Range("L114").Value = 47000
$ws.Range("N114").Value = -55678
$ws.Range("H116").Value = 2078.889
$ws.Range("I116").Value = 1393.3334
$ws.Range("J116").Value = 3450
$ws.Range("K116").Value = 1393.3334
$ws.Range("L116").Value = 3450
$ws.Range("M116").Value = 900.6666
$ws.Range("N116").Value = -8038
$ws.Range("H117").Value = 98596
$ws.Range("J117").Value = 98596
$ws.Range("L117").Value = 98596
$ws.Range("N117").Value = -107774
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 79982
$ws.Range("J120").Value = 79982
$ws.Range("L120").Value = 79982
$ws.Range("N120").Value = -89658
$ws.Range("H122").Value = 3675.3784
$ws.Range("I122").Value = 2428.3572
$ws.Range("J122").Value = 7555
$ws.Range("K122").Value = 7285.071599999999
$ws.Range("L122").Value = 22665
$ws.Range("M122").Value = -4835.071599999999
$ws.Range("N122").Value = -27565
$ws.Range("H123").Value = 101994.5
$ws.Range("J123").Value = 101994.5
$ws.Range("L123").Value = 101994.5
$ws.Range("N123").Value = -111794.5
$ws.Range("H125").Value = 114326
$ws.Range("J125").Value = 114326
$ws.Range("L125").Value = 114326
$ws.Range("N125").Value = -124166
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H128").Value = 139000
$ws.Range("J128").Value = 139000
$ws.Range("L128").Value = 139000
$ws.Range("N128").Value = -148960
$ws.Range("H129").Value = 141334
$ws.Range("J129").Value = 141334
$ws.Range("L129").Value = 141334
$ws.Range("N129").Value = -151334
$ws.Range("H136").Value = 3008.611
$ws.Range("I136").Value = 2697.25
$ws.Range("K136").Value = 8091.75
$ws.Range("M136").Value = -5541.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2078.889
$ws.Range("I3").Value = 1393.3334
$ws.Range("J3").Value = 3450
$ws.Range("K3").Value = 1393.3334
$ws.Range("L3").Value = 3450
$ws.Range("M3").Value = -1279.3334
$ws.Range("N3").Value = -3678
$ws.Range("H86").Value = 642.75
$ws.Range("I86").Value = 486.5
$ws.Range("K86").Value = 486.5
$ws.Range("M86").Value = 636.5
$ws.Range("H89").Value = 642.75
$ws.Range("I89").Value = 486.5
$ws.Range("K89").Value = 2432.5
$ws.Range("M89").Value = 3183.5
$ws.Range("H105").Value = 2042.55
$ws.Range("I105").Value = 1736.9412
$ws.Range("K105").Value = 1736.9412
$ws.Range("M105").Value = 10.05880000000002
$ws.Range("H134").Value = 18521336
$ws.Range("I134").Value = 2611.111
$ws.Range("K134").Value = 7833.333
$ws.Range("M134").Value = -5298.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 135.76923
$ws.Range("I7").Value = 68.454544
$ws.Range("K7").Value = 68.454544
$ws.Range("M7").Value = 44.545456
$ws.Range("H16").Value = 2281.25
$ws.Range("I16").Value = 1875
$ws.Range("K16").Value = 1875
$ws.Range("M16").Value = -1588
$ws.Range("H22").Value = 566.3333
$ws.Range("I22").Value = 449.5
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 449.5
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -99.5
$ws.Range("N22").Value = -1500
$ws.Range("H31").Value = 2812.963
$ws.Range("I31").Value = 2134.8635
$ws.Range("J31").Value = 5796.6
$ws.Range("K31").Value = 2134.8635
$ws.Range("L31").Value = 5796.6
$ws.Range("M31").Value = -1839.8635
$ws.Range("N31").Value = -6386.6
$ws.Range("H34").Value = 2812.963
$ws.Range("I34").Value = 2134.8635
$ws.Range("J34").Value = 5796.6
$ws.Range("K34").Value = 2134.8635
$ws.Range("L34").Value = 5796.6
$ws.Range("M34").Value = -1932.8635
$ws.Range("N34").Value = -6200.6
$ws.Range("H58").Value = 3133.84
$ws.Range("I58").Value = 2773.7715
$ws.Range("J58").Value = 3974
$ws.Range("K58").Value = 2773.7715
$ws.Range("L58").Value = 3974
$ws.Range("M58").Value = -2570.7715
$ws.Range("N58").Value = -4380
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H76").Value = 8208.143
$ws.Range("I76").Value = 8208.143
$ws.Range("K76").Value = 8208.143
$ws.Range("M76").Value = -7893.143
$ws.Range("H79").Value = 8208.143
$ws.Range("I79").Value = 8208.143
$ws.Range("K79").Value = 8208.143
$ws.Range("M79").Value = -7116.143
$ws.Range("H94").Value = 1937.9445
$ws.Range("I94").Value = 1248.5
$ws.Range("J94").Value = 2024.125
$ws.Range("K94").Value = 1248.5
$ws.Range("L94").Value = 2024.125
$ws.Range("M94").Value = -797.5
$ws.Range("N94").Value = -2926.125
$ws.Range("H113").Value = 2281.25
$ws.Range("I113").Value = 1875
$ws.Range("K113").Value = 1875
$ws.Range("M113").Value = 295
$ws.Range("H116").Value = 85329.336
$ws.Range("J116").Value = 85329.336
$ws.Range("L116").Value = 85329.336
$ws.Range("N116").Value = -94507.336
$ws.Range("H122").Value = 1481.1578
$ws.Range("I122").Value = 1582.8
$ws.Range("K122").Value = 4748.4
$ws.Range("M122").Value = -2298.4
$ws.Range("H134").Value = 2973.96
$ws.Range("I134").Value = 2373.8572
$ws.Range("K134").Value = 7121.571599999999
$ws.Range("M134").Value = -4586.571599999999
$ws.Range("H136").Value = 3133.84
$ws.Range("I136").Value = 2773.7715
$ws.Range("J136").Value = 3974
$ws.Range("K136").Value = 8321.3145
$ws.Range("L136").Value = 11922
$ws.Range("M136").Value = -5771.3145
$ws.Range("N136").Value = -17022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2478.3333
$ws.Range("I17").Value = 5125
$ws.Range("K17").Value = 15375
$ws.Range("M17").Value = -15206
$ws.Range("H58").Value = 9415.308000000001
$ws.Range("I58").Value = 9415.308000000001
$ws.Range("K58").Value = 28245.924
$ws.Range("M58").Value = -28117.924
$ws.Range("H68").Value = 1667.6666
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2003
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 6009
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -7631
$ws.Range("H71").Value = 1667.6666
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2003
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 18027
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -26139
$ws.Range("H99").Value = 1699.6666
$ws.Range("I99").Value = 1699.6666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5098.9998
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2852.9998
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 403.63635
$ws.Range("I122").Value = 402
$ws.Range("J122").Value = 403.8
$ws.Range("K122").Value = 3618
$ws.Range("L122").Value = 3634.2
$ws.Range("M122").Value = -1168
$ws.Range("N122").Value = -8534.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H44").Value = 50031
$ws.Range("J44").Value = 50031
$ws.Range("L44").Value = 50031
$ws.Range("N44").Value = -51223
$ws.Range("H70").Value = 4500
$ws.Range("I70").Value = 4500
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 4500
$ws.Range("I73").Value = 4500
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -6372
$ws.Range("H97").Value = 1008.5455
$ws.Range("I97").Value = 427.16666
$ws.Range("K97").Value = 427.16666
$ws.Range("M97").Value = 68.83334000000002
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378
$ws.Range("H107").Value = 200.05882
$ws.Range("I107").Value = 148.22223
$ws.Range("J107").Value = 258.375
$ws.Range("K107").Value = 148.22223
$ws.Range("L107").Value = 258.375
$ws.Range("M107").Value = 1771.77777
$ws.Range("N107").Value = -4098.375
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H114").Value = 79191.8
$ws.Range("J114").Value = 60739.75
$ws.Range("L114").Value = 60739.75
$ws.Range("N114").Value = -69417.75
$ws.Range("H122").Value = 1760.7826
$ws.Range("I122").Value = 1389
$ws.Range("J122").Value = 1999.7858
$ws.Range("K122").Value = 4167
$ws.Range("L122").Value = 5999.357400000001
$ws.Range("M122").Value = -1717
$ws.Range("N122").Value = -10899.3574

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2934.2354
$ws.Range("I22").Value = 2088.5881
$ws.Range("K22").Value = 2088.5881
$ws.Range("M22").Value = -1793.5881
$ws.Range("H27").Value = 2934.2354
$ws.Range("I27").Value = 2088.5881
$ws.Range("K27").Value = 2088.5881
$ws.Range("M27").Value = -1981.5881
$ws.Range("H38").Value = 25015
$ws.Range("I38").Value = 20030
$ws.Range("K38").Value = 20030
$ws.Range("M38").Value = -19620
$ws.Range("H40").Value = 4496.615
$ws.Range("I40").Value = 4410.778
$ws.Range("J40").Value = 4689.75
$ws.Range("K40").Value = 4410.778
$ws.Range("L40").Value = 4689.75
$ws.Range("M40").Value = -4274.778
$ws.Range("N40").Value = -4961.75
$ws.Range("H46").Value = 2756.4167
$ws.Range("J46").Value = 2938.818
$ws.Range("L46").Value = 2938.818
$ws.Range("N46").Value = -3314.818
$ws.Range("H55").Value = 267.85294
$ws.Range("I55").Value = 228.75
$ws.Range("J55").Value = 302.6111
$ws.Range("K55").Value = 228.75
$ws.Range("L55").Value = 302.6111
$ws.Range("M55").Value = -55.75
$ws.Range("N55").Value = -648.6111000000001
$ws.Range("H68").Value = 3349.6
$ws.Range("I68").Value = 3285.1428
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 3285.1428
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2536.1428
$ws.Range("N68").Value = -4998
$ws.Range("H71").Value = 3349.6
$ws.Range("I71").Value = 3285.1428
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 16425.714
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -12681.714
$ws.Range("N71").Value = -24988
$ws.Range("H122").Value = 26655.652
$ws.Range("I122").Value = 22150.525
$ws.Range("J122").Value = 48055
$ws.Range("K122").Value = 66451.57500000001
$ws.Range("L122").Value = 144165
$ws.Range("M122").Value = -64001.57500000001
$ws.Range("N122").Value = -149065
$ws.Range("H132").Value = 3385.4
$ws.Range("I132").Value = 2909.1667
$ws.Range("K132").Value = 8727.500100000001
$ws.Range("M132").Value = -6197.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 69996
$ws.Range("J110").Value = 69996
$ws.Range("L110").Value = 69996
$ws.Range("N110").Value = -78176
$ws.Range("H122").Value = 4391.25
$ws.Range("I122").Value = 1876.7142
$ws.Range("J122").Value = 10258.5
$ws.Range("K122").Value = 5630.142599999999
$ws.Range("L122").Value = 30775.5
$ws.Range("M122").Value = -3180.142599999999
$ws.Range("N122").Value = -35675.5
$ws.Range("H126").Value = 39992
$ws.Range("I126").Value = 39992
$ws.Range("K126").Value = 119976
$ws.Range("M126").Value = -117506
$ws.Range("H136").Value = 1688.0193
$ws.Range("I136").Value = 1085.9756
$ws.Range("J136").Value = 3932
$ws.Range("K136").Value = 3257.9268
$ws.Range("L136").Value = 11796
$ws.Range("M136").Value = -707.9268000000002
$ws.Range("N136").Value = -16896
